# Generate Report for Handoff
#
# This mirrors a "generate handoff report" run: the previous handback
# ("Handed back: in sync with en-US" @ 2016-08-25 20:46:53) is superseded by
# a fresh "Ready for handoff" status (@ 2016-08-25 20:47:55) on the Overview
# sheet, and the per-locale sheets record that the most recent handback file
# for the second (bef00e20-...) document is now considered stale, with a new
# "Latest Handback DateTime" and populated "Error Detail" describing the
# mismatch.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 (bef00e20 document) now shows "Ready for handoff" ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-25 20:47:55"

# --- zh-cn sheet: row 3 (bef00e20 document) ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-25 20:47:51"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a151a12ed581cc9b25c6cfe0b8b47d48165c796/e2e/bef00e20-43ab-4ba6-8eee-8cfbff941922.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e931a6a2c17453d60768b82922722c069362f10/e2e/bef00e20-43ab-4ba6-8eee-8cfbff941922.md."

# --- de-de sheet: row 3 (bef00e20 document) ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-25 20:47:55"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a151a12ed581cc9b25c6cfe0b8b47d48165c796/e2e/bef00e20-43ab-4ba6-8eee-8cfbff941922.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e931a6a2c17453d60768b82922722c069362f10/e2e/bef00e20-43ab-4ba6-8eee-8cfbff941922.md."

# --- Widen "Error Detail" column (P) on both locale sheets so the new long
#     message is readable. NOTE: Excel's ColumnWidth property (in "characters")
#     is not a 1:1 match with the stored OOXML column width (which includes
#     fixed pixel padding); 39.1666... is the COM value that round-trips to a
#     stored width of exactly 40, matching the other width="40" columns in
#     this sheet (e.g. columns A, G, I). ---
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666667
$dede.Columns.Item(16).ColumnWidth = 39.1666666667
